$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style/formatting from an existing header cell (AC1) onto
# the new header cells so they match (bold, centered, bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record for every data row (2-46): Wins=91, Losses=71, Ties=0
$ws.Range("AD2:AD46").Value = 91
$ws.Range("AE2:AE46").Value = 71
$ws.Range("AF2:AF46").Value = 0
